$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "Carlo Pomarolli"
$ws.Range("B45").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C45").Value = "Luca Frasca | Clitoriders"
$ws.Range("D45").Value = "Randy Cobbinah | MAI UNA GIOIA"
$ws.Range("E45").Value = "Michele Parisi  | MediaserT"
$ws.Range("F45").Value = "Davide  Bazzano  | iMontagna"
